# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt -
# Albahaca" right after the existing row 60 (new row 61), pushing the old
# rows 61-94 down to 62-95.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61; existing row 61 (and everything below)
# shifts down by one, carrying its formatting along (Excel's default
# "insert" behaviour).
$ws.Rows("61:61").Insert()

# Populate the newly inserted row 61 with the new weekly record.
$ws.Range("A61").Value = 4
$ws.Range("B61").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C61").Value = "Los Lagos"
$ws.Range("D61").Value = 44596
$ws.Range("E61").Value = 10
$ws.Range("F61").Value = 100112052
$ws.Range("G61").Value = "Albahaca"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 120
$ws.Range("K61").Value = 6000
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = 6000
$ws.Range("N61").Value = "$/docena de matas"
$ws.Range("O61").Value = "Región Metropolitana"
$ws.Range("P61").Value = 1000
$ws.Range("Q61").Value = 6
$ws.Range("R61").Value = "Hortaliza"
